$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D to hold the "Status" field
$ws.Columns.Item(4).Insert()

# Update header row (row 1)
$headerRow = @("ISIN", "Stock Name", "Mutual Fund", "Status", "Jan_2026", "Dec_2025", "Oct_2025", "MoM", "QoQ")
for ($i = 0; $i -lt $headerRow.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headerRow[$i]
}

# Data rows
$data = @(
    @("INE931S01010", "Adani Energy Solutions Limited", "quant Commodities Fund", "Reducing", 9.512864, 9.812553, 2.438588, -0.299688999999999, 7.074276),
    @("INE423A01024", "Adani Enterprises Limited", "quant Commodities Fund", "Adding Consistently", 9.024507, 8.987818, 1.919013, 0.03668899999999908, 7.105494),
    @("INE371A01025", "Graphite India Limited", "quant Commodities Fund", "Reducing", 8.991999, 10.366244, 0.0, -1.374245, 8.991999),
    @("INE364U01010", "Adani Green Energy Limited", "quant Commodities Fund", "Reducing Consistently", 8.815982, 8.948149, 9.409195, -0.1321670000000008, -0.5932130000000004),
    @("INE281B01032", "Lloyds Metals And Energy Limited", "quant Commodities Fund", "Reducing", 8.211117, 8.830628, 5.482282, -0.619511000000001, 2.728835),
    @("INE079A01024", "Ambuja Cements Ltd", "quant Commodities Fund", "Adding Consistently", 5.6717, 5.557017, 0.0, 0.1146830000000003, 5.6717),
    @("INE331A01037", "The Ramco Cements Limited", "quant Commodities Fund", "Fresh Entry", 4.791299, 0.0, 0.0, 4.791299, 4.791299),
    @("INE907A01026", "Kalyani Steels Ltd", "quant Commodities Fund", "Adding", 4.422661, 4.310524, 4.530184, 0.1121369999999997, -0.1075230000000005),
    @("INE081A01020", "Tata Steel Limited", "quant Commodities Fund", "Fresh Entry", 3.04389, 0.0, 0.0, 3.04389, 3.04389),
    @("INE628A01036", "UPL Limited", "quant Commodities Fund", "Reducing", 2.922882, 2.966176, 0.0, -0.04329399999999994, 2.922882),
    @("INE090A01021", "ICICI Bank Limited", "quant Commodities Fund", "Fresh Entry", 2.717233, 0.0, 0.0, 2.717233, 2.717233),
    @("INE373A01013", "BASF India Ltd", "quant Commodities Fund", "Complete Exit", 0.0, 0.0, 2.97225, 0.0, -2.97225),
    @("INE814H01029", "Adani Power Limited", "quant Commodities Fund", "Complete Exit", 0.0, 2.097552, 9.571447, -2.097552, -9.571447),
    @("INE813H01021", "Torrent Power Limited", "quant Commodities Fund", "Complete Exit", 0.0, 1.246209, 0.0, -1.246209, 0.0),
    @("INE752E01010", "Power Grid Corporation of India Limited", "quant Commodities Fund", "Complete Exit", 0.0, 0.0, 1.002811, 0.0, -1.002811),
    @("INE462A01022", "Bayer Cropscience Ltd", "quant Commodities Fund", "Complete Exit", 0.0, 0.0, 7.250817, 0.0, -7.250817),
    @("INE075A01022", "Wipro Ltd", "quant Commodities Fund", "Complete Exit", 0.0, 1.816168, 0.0, -1.816168, 0.0),
    @("INE267A01025", "Hindustan Zinc Limited", "quant Commodities Fund", "Complete Exit", 0.0, 6.836897, 0.0, -6.836897, 0.0),
    @("INE258G01013", "Sumitomo Chemical India Limited", "quant Commodities Fund", "Complete Exit", 0.0, 0.0, 4.134973, 0.0, -4.134973),
    @("INE257A01026", "Bharat Heavy Electricals Ltd", "quant Commodities Fund", "Complete Exit", 0.0, 2.903893, 0.0, -2.903893, 0.0),
    @("INE245A01021", "Tata Power Company Limited", "quant Commodities Fund", "Complete Exit", 0.0, 1.895742, 9.665714, -1.895742, -9.665714),
    @("INE213A01029", "Oil and Natural Gas Corporation Ltd.", "quant Commodities Fund", "Complete Exit", 0.0, 0.0, 7.423331, 0.0, -7.423331),
    @("INE200A01026", "GE Vernova T&D India Limited", "quant Commodities Fund", "Complete Exit", 0.0, 2.837839, 0.0, -2.837839, 0.0),
    @("INE129A01019", "GAIL (India) Limited", "quant Commodities Fund", "Complete Exit", 0.0, 0.0, 6.301544, 0.0, -6.301544),
    @("INE296A01032", "Bajaj Finance Limited", "quant Commodities Fund", "Complete Exit", 0.0, 0.0, 8.812978, 0.0, -8.812978),
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
